$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12901
$ws1.Range("F3").Value = 634
$ws1.Range("F6").Value = 333
$ws1.Range("F7").Value = 411
$ws1.Range("F8").Value = 238
$ws1.Range("F9").Value = 13064
$ws1.Range("F10").Value = 47
$ws1.Range("F11").Value = 38
$ws1.Range("F12").Value = 5302
$ws1.Range("F16").Value = 36
$ws1.Range("F17").Value = 1205
$ws1.Range("F20").Value = 689
$ws1.Range("F21").Value = 2866
$ws1.Range("F22").Value = 6232
$ws1.Range("F23").Value = 1168
$ws1.Range("F24").Value = 3643

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12901
$ws4.Range("F3").Value = 634
$ws4.Range("F6").Value = 333
$ws4.Range("F8").Value = 411
$ws4.Range("F9").Value = 238
$ws4.Range("F10").Value = 13064
$ws4.Range("F11").Value = 47
$ws4.Range("F12").Value = 38
$ws4.Range("F13").Value = 5302
$ws4.Range("F17").Value = 36
$ws4.Range("F18").Value = 1205
$ws4.Range("F21").Value = 689
$ws4.Range("F22").Value = 2866
$ws4.Range("F24").Value = 6232
$ws4.Range("F25").Value = 1168
$ws4.Range("F26").Value = 3643
